$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(2, 2).Value = 24038906.3850759
$ws.Cells.Item(3, 2).Value = 28931.17366267106
$ws.Cells.Item(4, 2).Value = 29370.96996600106
$ws.Cells.Item(5, 2).Value = 29820.12961069287
$ws.Cells.Item(6, 2).Value = 30278.97595935772
$ws.Cells.Item(7, 2).Value = 30747.84532008644
$ws.Cells.Item(8, 2).Value = 31227.08730381621
$ws.Cells.Item(9, 2).Value = 31717.06517595852
$ws.Cells.Item(10, 2).Value = 32218.15620086791
$ws.Cells.Item(11, 2).Value = 32730.75197723418
$ws.Cells.Item(12, 2).Value = 33255.25876313686
$ws.Cells.Item(13, 2).Value = 33792.09778907381
$ws.Cells.Item(14, 2).Value = 34341.70555676104
$ws.Cells.Item(15, 2).Value = 34904.53412245605
$ws.Cells.Item(16, 2).Value = 35481.0513622764
$ws.Cells.Item(17, 2).Value = 36071.74121807304
$ws.Cells.Item(18, 2).Value = 36677.10392130951
$ws.Cells.Item(19, 2).Value = 37297.656193139
$ws.Cells.Item(20, 2).Value = 37933.93141843074
$ws.Cells.Item(21, 2).Value = 38586.47979127956
$ws.Cells.Item(22, 2).Value = 39255.86843015016
$ws.Cells.Item(23, 2).Value = 39942.68145983873
$ws.Cells.Item(24, 2).Value = 40647.52005818983
$ws.Cells.Item(25, 2).Value = 41371.00246518651
$ws.Cells.Item(26, 2).Value = 42113.76395169138
$ws.Cells.Item(27, 2).Value = 42876.45674582875
$ws.Cells.Item(28, 2).Value = 43659.74991421455
$ws.Cells.Item(29, 2).Value = 44464.32919577739
$ws.Cells.Item(30, 2).Value = 45290.8967858418
$ws.Cells.Item(31, 2).Value = 46140.17106781359
$ws.Cells.Item(32, 2).Value = 47012.88629040297
$ws.Cells.Item(33, 2).Value = 47909.79218796983
$ws.Cells.Item(34, 2).Value = 48831.65354171403
$ws.Cells.Item(35, 2).Value = 49779.24967969092
$ws.Cells.Item(36, 2).Value = 50753.37391341887
$ws.Cells.Item(37, 2).Value = 51754.83290932302
$ws.Cells.Item(38, 2).Value = 52784.44599287015
$ws.Cells.Item(39, 2).Value = 53843.04438405027
$ws.Cells.Item(40, 2).Value = 54931.47036235551
$ws.Cells.Item(41, 2).Value = 56050.57635997441
$ws.Cells.Item(42, 2).Value = 57201.22398179524
$ws.Cells.Item(43, 2).Value = 58384.28295163765
$ws.Cells.Item(44, 2).Value = 59600.62998317974
$ws.Cells.Item(45, 2).Value = 60851.14757564523
$ws.Cells.Item(46, 2).Value = 62136.72273341841
$ws.Cells.Item(47, 2).Value = 63458.24560953621
$ws.Cells.Item(48, 2).Value = 64816.60807353072
$ws.Cells.Item(49, 2).Value = 66212.70220344847
$ws.Cells.Item(50, 2).Value = 67647.41870317706
$ws.Cells.Item(51, 2).Value = 69121.64524604419
$ws.Cells.Item(52, 2).Value = 70636.2647457352
$ws.Cells.Item(53, 2).Value = 72192.15355647876
$ws.Cells.Item(54, 2).Value = 73790.17960416783
$ws.Cells.Item(55, 2).Value = 75431.20045115406
$ws.Cells.Item(56, 2).Value = 77116.06129681887
$ws.Cells.Item(57, 2).Value = 78845.59291737618
$ws.Cells.Item(58, 2).Value = 80620.60954829578
$ws.Cells.Item(59, 2).Value = 82441.90671264996
$ws.Cells.Item(60, 2).Value = 84310.25899999785
$ws.Cells.Item(61, 2).Value = 86226.41780011912
$ws.Cells.Item(62, 2).Value = 88191.10899634421
$ws.Cells.Item(63, 2).Value = 90205.03062373254
$ws.Cells.Item(64, 2).Value = 92268.85049802369
$ws.Cells.Item(65, 2).Value = 94383.20382054512
$ws.Cells.Item(66, 2).Value = 96548.69076636239
$ws.Cells.Item(67, 2).Value = 98765.8740612913
$ws.Cells.Item(68, 2).Value = 101035.2765554267
$ws.Cells.Item(69, 2).Value = 103357.37880023
$ws.Cells.Item(70, 2).Value = 105732.616636596
$ws.Cells.Item(71, 2).Value = 108161.3788017988
$ws.Cells.Item(72, 2).Value = 110644.0045636592
$ws.Cells.Item(73, 2).Value = 113180.7813897605
$ws.Cells.Item(74, 2).Value = 115771.9426607393
$ws.Cells.Item(75, 2).Value = 118417.6654360648
$ws.Cells.Item(76, 2).Value = 121118.0682812925
$ws.Cells.Item(77, 2).Value = 123873.2091657631
$ws.Cells.Item(78, 2).Value = 126683.0834399593
$ws.Cells.Item(79, 2).Value = 129547.6219014596
$ws.Cells.Item(80, 2).Value = 132466.6889589678
$ws.Cells.Item(81, 2).Value = 135440.0809033826
$ws.Cells.Item(82, 2).Value = 138467.524295167
$ws.Cells.Item(83, 2).Value = 141548.6744770557
$ws.Cells.Item(84, 2).Value = 144683.1142210836
$ws.Cells.Item(85, 2).Value = 147870.3525186485
$ws.Cells.Item(86, 2).Value = 151109.8235223399
$ws.Cells.Item(87, 2).Value = 154400.885647508
$ws.Cells.Item(88, 2).Value = 157742.8208422041
$ws.Cells.Item(89, 2).Value = 161134.8340328744
$ws.Cells.Item(90, 2).Value = 164576.0527528782
$ws.Cells.Item(91, 2).Value = 168065.5269615471
$ws.Cells.Item(92, 2).Value = 171602.2290596815
$ws.Cells.Item(93, 2).Value = 175185.0541075802
$ws.Cells.Item(94, 2).Value = 178812.8202513684
$ws.Cells.Item(95, 2).Value = 182484.2693623695
$ws.Cells.Item(96, 2).Value = 186198.0678936799
$ws.Cells.Item(97, 2).Value = 189952.8079579817
$ws.Cells.Item(98, 2).Value = 193747.0086293317
$ws.Cells.Item(99, 2).Value = 197579.1174713788
$ws.Cells.Item(100, 2).Value = 201447.5122935555
$ws.Cells.Item(101, 2).Value = 205350.5031359127
$ws.Cells.Item(102, 2).Value = 209286.3344831245
$ws.Cells.Item(103, 2).Value = 213253.1877063259
$ws.Cells.Item(104, 2).Value = 217249.183731465
$ws.Cells.Item(105, 2).Value = 221272.3859321581
$ws.Cells.Item(106, 2).Value = 225320.8032429982
$ws.Cells.Item(107, 2).Value = 229392.3934901218
$ws.Cells.Item(108, 2).Value = 233485.0669333382
$ws.Cells.Item(109, 2).Value = 237596.6900146205
$ws.Cells.Item(110, 2).Value = 241725.0893060436
$ws.Cells.Item(111, 2).Value = 245868.055649659
$ws.Cells.Item(112, 2).Value = 250023.3484812781
$ws.Cells.Item(113, 2).Value = 254188.7003286669
$ws.Cells.Item(114, 2).Value = 258361.8214747118
$ws.Cells.Item(115, 2).Value = 262540.4047744976
$ws.Cells.Item(116, 2).Value = 266722.1306150863
$ws.Cells.Item(117, 2).Value = 270904.672005857
$ws.Cells.Item(118, 2).Value = 275085.6997866787
$ws.Cells.Item(119, 2).Value = 279262.8879402653
$ws.Cells.Item(120, 2).Value = 283433.9189950014
$ws.Cells.Item(121, 2).Value = 287596.4895034602
$ws.Cells.Item(122, 2).Value = 291748.3155816759
$ws.Cells.Item(123, 2).Value = 295887.1384935569
$ws.Cells.Item(124, 2).Value = 300010.7302646402
$ws.Cells.Item(125, 2).Value = 304116.8993089068
$ws.Cells.Item(126, 2).Value = 308203.496052207
$ws.Cells.Item(127, 2).Value = 312268.4185354366
$ws.Cells.Item(128, 2).Value = 316309.6179808676
$ws.Cells.Item(129, 2).Value = 320325.1043042946
$ws.Cells.Item(130, 2).Value = 324312.9515565767
$ws.Cells.Item(131, 2).Value = 328271.3032769083
$ws.Cells.Item(132, 2).Value = 332198.3777418604
$ws.Cells.Item(133, 2).Value = 336092.473093009
$ws.Cells.Item(134, 2).Value = 339951.9723267516
$ws.Cells.Item(135, 2).Value = 343775.348130982
$ws.Cells.Item(136, 2).Value = 347561.1675518057
$ws.Cells.Item(137, 2).Value = 351308.0964763212
$ws.Cells.Item(138, 2).Value = 355014.9039159936
$ws.Cells.Item(139, 2).Value = 358680.4660766122
$ws.Cells.Item(140, 2).Value = 362303.7702018087
$ws.Cells.Item(141, 2).Value = 365883.9181768732
$ws.Cells.Item(142, 2).Value = 369420.129881015
$ws.Cells.Item(143, 2).Value = 372911.7462771733
$ws.Cells.Item(144, 2).Value = 376358.2322284224
$ws.Cells.Item(145, 2).Value = 379759.1790322061
$ws.Cells.Item(146, 2).Value = 383114.3066631944
$ws.Cells.Item(147, 2).Value = 386423.4657180435
$ws.Cells.Item(148, 2).Value = 389686.6390546193
$ws.Cells.Item(149, 2).Value = 392903.9431210484
$ws.Cells.Item(150, 2).Value = 396075.6289698171
$ws.Cells.Item(151, 2).Value = 399202.0829533363
$ws.Cells.Item(152, 2).Value = 402283.8270994871
$ws.Cells.Item(153, 2).Value = 405321.5191652137
$ws.Cells.Item(154, 2).Value = 408315.9523689653
$ws.Cells.Item(155, 2).Value = 411268.0548024827
$ws.Cells.Item(156, 2).Value = 414178.8885244788
$ws.Cells.Item(157, 2).Value = 417049.64833991
$ws.Cells.Item(158, 2).Value = 419881.6602683924
$ws.Cells.Item(159, 2).Value = 422676.3797084171
$ws.Cells.Item(160, 2).Value = 425435.3893030105
$ws.Cells.Item(161, 2).Value = 428160.3965150334
$ws.Cells.Item(162, 2).Value = 430853.230920457
$ws.Cells.Item(163, 2).Value = 433515.841229245
$ws.Cells.Item(164, 2).Value = 436150.2920444395
$ws.Cells.Item(165, 2).Value = 438758.7603703769
$ws.Cells.Item(166, 2).Value = 441343.53188277
$ws.Cells.Item(167, 2).Value = 443906.996973303
$ws.Cells.Item(168, 2).Value = 446451.6465822417
$ws.Cells.Item(169, 2).Value = 448980.0678336342
$ws.Cells.Item(170, 2).Value = 451494.9394880468
$ws.Cells.Item(171, 2).Value = 453999.0272280248
$ws.Cells.Item(172, 2).Value = 456495.1787924619
$ws.Cells.Item(173, 2).Value = 458986.3189761405
$ws.Cells.Item(174, 2).Value = 461475.4445111806
$ws.Cells.Item(175, 2).Value = 463965.6188472253
$ws.Cells.Item(176, 2).Value = 466459.9668477707
$ws.Cells.Item(177, 2).Value = 468961.6694198764
$ws.Cells.Item(178, 2).Value = 471473.95809457
$ws.Cells.Item(179, 2).Value = 474000.1095755563
$ws.Cells.Item(180, 2).Value = 476543.4402733754
$ws.Cells.Item(181, 2).Value = 479107.3008424806
$ws.Cells.Item(182, 2).Value = 481695.0707381996
$ws.Cells.Item(183, 2).Value = 484310.1528103013
$ws.Cells.Item(184, 2).Value = 486955.967949711
$ws.Cells.Item(185, 2).Value = 489635.9498047909
$ws.Cells.Item(186, 2).Value = 492353.5395824419
$ws.Cells.Item(187, 2).Value = 495112.1809495632
$ws.Cells.Item(188, 2).Value = 497915.3150497563
$ws.Cells.Item(189, 2).Value = 500766.3756488123
$ws.Cells.Item(190, 2).Value = 503668.7844234371
$ws.Cells.Item(191, 2).Value = 506625.9464052642
$ws.Cells.Item(192, 2).Value = 509641.2455926861
$ws.Cells.Item(193, 2).Value = 512718.0407424713
$ws.Cells.Item(194, 2).Value = 515859.6613508855
$ws.Cells.Item(195, 2).Value = 519069.4038350928
$ws.Cells.Item(196, 2).Value = 522350.5279238387
$ws.Cells.Item(197, 2).Value = 525706.2532654905
$ws.Cells.Item(198, 2).Value = 529139.7562611412
$ws.Cells.Item(199, 2).Value = 532654.1671295586
$ws.Cells.Item(200, 2).Value = 536252.5672094484
$ws.Cells.Item(201, 2).Value = 539937.9865044794
$ws.Cells.Item(202, 2).Value = 543713.4014748693
$ws.Cells.Item(203, 2).Value = 547581.8552689541
$ws.Cells.Item(204, 2).Value = 551545.9806156873
$ws.Cells.Item(205, 2).Value = 555608.692443642
$ws.Cells.Item(206, 2).Value = 559772.7360446659
$ws.Cells.Item(207, 2).Value = 564040.7959135106
$ws.Cells.Item(208, 2).Value = 568415.4950594479
$ws.Cells.Item(209, 2).Value = 572899.3945706685
$ws.Cells.Item(210, 2).Value = 577494.9934289779
$ws.Cells.Item(211, 2).Value = 582204.7285719563
$ws.Cells.Item(212, 2).Value = 587030.9751990071
$ws.Cells.Item(213, 2).Value = 591976.0473168263
$ws.Cells.Item(214, 2).Value = 597042.1985196359
$ws.Cells.Item(215, 2).Value = 602231.6229989832
$ws.Cells.Item(216, 2).Value = 607546.4567765237
$ws.Cells.Item(217, 2).Value = 612988.7791538185
$ws.Cells.Item(218, 2).Value = 618560.614372163
$ws.Cells.Item(219, 2).Value = 624263.9334747227
$ws.Cells.Item(220, 2).Value = 630100.6563635111
$ws.Cells.Item(221, 2).Value = 636072.6540435624
$ws.Cells.Item(222, 2).Value = 642181.7510440808
$ws.Cells.Item(223, 2).Value = 648429.7280104604
$ws.Cells.Item(224, 2).Value = 654818.3244555739
$ws.Cells.Item(225, 2).Value = 661349.2416634716
$ws.Cells.Item(226, 2).Value = 668024.1457334403
$ws.Cells.Item(227, 2).Value = 674844.6707581955
$ws.Cells.Item(228, 2).Value = 681812.4221231932
$ws.Cells.Item(229, 2).Value = 688928.9799196895
$ws.Cells.Item(230, 2).Value = 696195.9024605139
$ws.Cells.Item(231, 2).Value = 703614.729889658
$ws.Cells.Item(232, 2).Value = 711186.9878743537
$ws.Cells.Item(233, 2).Value = 718914.1913713531
$ws.Cells.Item(234, 2).Value = 726797.8484570335
$ws.Cells.Item(235, 2).Value = 734839.4642108064
$ws.Cells.Item(236, 2).Value = 743040.5446435378
$ws.Cells.Item(237, 2).Value = 751402.6006602478
$ws.Cells.Item(238, 2).Value = 759927.1520479799
$ws.Cells.Item(239, 2).Value = 768615.7314797018
$ws.Cells.Item(240, 2).Value = 777469.8885246343
$ws.Cells.Item(241, 2).Value = 786491.1936555288
$ws.Cells.Item(242, 2).Value = 795681.2422450188
$ws.Cells.Item(243, 2).Value = 805041.6585409378
$ws.Cells.Item(244, 2).Value = 814574.0996121336
$ws.Cells.Item(245, 2).Value = 824280.2592570143
$ws.Cells.Item(246, 2).Value = 834161.8718648576
$ws.Cells.Item(247, 2).Value = 844220.7162230125
$ws.Cells.Item(248, 2).Value = 854458.6192605492
$ws.Cells.Item(249, 2).Value = 864877.4597208505
$ws.Cells.Item(250, 2).Value = 875479.1717548643
$ws.Cells.Item(251, 2).Value = 886265.7484273619
$ws.Cells.Item(252, 2).Value = 897239.2451274528
$ws.Cells.Item(253, 2).Value = 908401.7828768232
$ws.Cells.Item(254, 2).Value = 919755.5515270047
$ws.Cells.Item(255, 2).Value = 931302.8128376763
$ws.Cells.Item(256, 2).Value = 943045.9034293685
$ws.Cells.Item(257, 2).Value = 954987.2376020161
$ws.Cells.Item(258, 2).Value = 967129.3100115269
$ws.Cells.Item(259, 2).Value = 979474.6981969101
$ws.Cells.Item(260, 2).Value = 992026.0649505258
$ws.Cells.Item(261, 2).Value = 1004786.160522058
$ws.Cells.Item(262, 2).Value = 1017757.824650417
$ws.Cells.Item(263, 2).Value = 1030943.98841326
$ws.Cells.Item(264, 2).Value = 1044347.675886897
$ws.Cells.Item(265, 2).Value = 1057972.005609593
$ws.Cells.Item(266, 2).Value = 1071820.191835804
$ws.Cells.Item(267, 2).Value = 1085895.545577675
$ws.Cells.Item(268, 2).Value = 1100201.475421216
$ws.Cells.Item(269, 2).Value = 1114741.488109022
$ws.Cells.Item(270, 2).Value = 1129519.188881505
$ws.Cells.Item(271, 2).Value = 1144538.281565264
$ws.Cells.Item(272, 2).Value = 1159802.568400727
$ws.Cells.Item(273, 2).Value = 1175315.949597442
$ws.Cells.Item(274, 2).Value = 1191082.422608247
$ws.Cells.Item(275, 2).Value = 1207106.081111367
$ws.Cells.Item(276, 2).Value = 1223391.113689453
$ws.Cells.Item(277, 2).Value = 1239941.802195994
$ws.Cells.Item(278, 2).Value = 1256762.519796714
$ws.Cells.Item(279, 2).Value = 1273857.728675801
$ws.Cells.Item(280, 2).Value = 1291231.977394898
$ws.Cells.Item(281, 2).Value = 1308889.89789407
$ws.Cells.Item(282, 2).Value = 1326836.202122384
$ws.Cells.Item(283, 2).Value = 1345075.678286856
$ws.Cells.Item(284, 2).Value = 1363613.186707621
$ws.Cells.Item(285, 2).Value = 1382453.655267628
$ws.Cells.Item(286, 2).Value = 1401602.074445377
$ws.Cells.Item(287, 2).Value = 1421063.491917483
$ws.Cells.Item(288, 2).Value = 1440843.006721368
$ws.Cells.Item(289, 2).Value = 1460945.762964756
$ws.Cells.Item(290, 2).Value = 1481376.943071615
$ws.Cells.Item(291, 2).Value = 1502141.760553045
$ws.Cells.Item(292, 2).Value = 1523245.452292206
$ws.Cells.Item(293, 2).Value = 1544693.270333078
$ws.Cells.Item(294, 2).Value = 1566490.473164077
$ws.Cells.Item(295, 2).Value = 1588642.316484787
$ws.Cells.Item(296, 2).Value = 1611154.043450316
$ws.Cells.Item(297, 2).Value = 1634030.874383345
$ws.Cells.Item(298, 2).Value = 1657277.995947354
$ws.Cells.Item(299, 2).Value = 1680900.549775077
$ws.Cells.Item(300, 2).Value = 1704903.620548153
$ws.Cells.Item(301, 2).Value = 1729292.223522967
$ws.Cells.Item(302, 2).Value = 1754071.291500279
